$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.1196539058413437
$ws.Range("B3").Value = 0.834450129110025
$ws.Range("B4").Value = 2257.292800218992
$ws.Range("B5").Value = '[0.9999999999999998, 0.203938050835155, 0.14038516634682768, 0.1106271597198941, 0.11919130928281732, 0.3652757667381568, 0.13534100875315108, 0.13130863884504743, 0.08710727075881138, 0.08630647578632664, 0.008921163593526227, 0.024739708966519464, 0.0233182526464828, 0.08361185834748502, 0.09106036177391073, 0.012178098189479762, 0.04188091879822754, 0.006598800571202339, 0.054506754862295316, 0.053717727415923185]'

$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 0.9847826086956522
$ws.Range("B4").Value = -17.43605063666466
$ws.Range("B5").Value = 0.7214705347193457
$ws.Range("B6").Value = 1544.408806512228
$ws.Range("B7").Value = 13.62227738077903
$ws.Range("B8").Value = 0.7704497792522496
$ws.Range("B9").Value = 1725.850087644882
$ws.Range("B10").Value = '[1.0, 0.23470523125727913, 0.18043162846783606, 0.13064911174831265, 0.1537980327107822, 0.2747986897706073, 0.10781826043880667, 0.10913547544620399, 0.0738815341274262, 0.0990760412509017, 0.026031972986113202, 0.03437277013119131, 0.01955173719188702, 0.054342646834331904, 0.06869620746520537, 0.016322515897666468, 0.045699809381399524, 0.011316028327488913, 0.04927894790543823, 0.026943033982928816]'

$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.2107369527606112
$ws.Range("B3").Value = 136.0711203183615
$ws.Range("B4").Value = 0.09616636969425546
$ws.Range("B5").Value = 0.8487813673533474
$ws.Range("B6").Value = '[1.0, 0.00138922602860478, -0.010978747279480927, -0.032401262472769345, -0.046392324809179236, 0.33919987925995787, 0.010851562939929624, 0.029608335052882368, -0.013131847850927799, -0.022627153426887748, -0.058741691954163996, -0.046364007036134997, -0.018437022134476545, 0.02263470621597901, 0.04356898547795632, -0.016637546046719586, 0.005802636495733714, -0.04105934097087084, 0.021714737130916824, 0.01530343816820269]'

$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.924108750216508
$ws.Range("B3").Value = 132.2180241842929
$ws.Range("B4").Value = 0.1222394146688099
$ws.Range("B5").Value = 0.8509908505480588
$ws.Range("B7").Value = -0.05170053783843203
$ws.Range("B6").Value = '[1.0, 0.0023056595849726036, -0.014430689943621384, -0.0358055588059253, -0.04517725009095377, 0.33813610618770995, 0.012313524922644306, 0.031197642051828404, -0.013159853039356545, -0.021417828132655782, -0.058997388905545206, -0.04554416692476793, -0.01771531533576334, 0.022008683470226144, 0.045657916989898795, -0.01798623415611631, 0.0063856828158397215, -0.04127922403191442, 0.021287687376868727, 0.015688677609790338]'

$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.54049807815235
$ws.Range("B3").Value = 163.0821755253245
$ws.Range("B4").Value = 0.1337560541901572
$ws.Range("B5").Value = 0.8191686387534715
$ws.Range("B7").Value = -0.05199194421438309
$ws.Range("B9").Value = 0.8159549426249337
$ws.Range("B6").Value = '[1.0, 0.03574753065209955, -0.007970617824398377, -0.007651802674355969, -0.024210834973450116, 0.256996855719331, 0.008379222938175673, 0.027502182818854026, -0.011896147276176352, 0.0074035776977858184, -0.048786900393601694, -0.040960863519489236, -0.039083975303286764, 0.02391688054247091, 0.04999596458174584, -0.029213394057002063, 0.0043563103872698815, -0.02794478656204016, 0.023629307708343508, 0.016003299296851712]'
